$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet3")

# Copy the existing header/data formatting onto the new columns before
# filling them in, so the new cells pick up the same styles (fonts,
# wrap text, etc.) as the existing A:C columns instead of creating new
# style/font table entries.
$ws.Range("A1:C2").Copy()
$ws.Range("D1").PasteSpecial(-4122)

# New header row (D1:F1)
$ws.Range("D1").Value = "exDay2"
$ws.Range("E1").Value = "exMonth2"
$ws.Range("F1").Value = "exYear2"

# New data row (D2:F2)
$ws.Range("D2").Value = 9
$ws.Range("E2").Value = "JUL"
$ws.Range("F2").Value = 2022

# Match the saved selection state in the target workbook.
$ws.Range("E2").Select()
